$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Row 29
$ws.Range("H29").Value = 1080
$ws.Range("I29").Value = 850
$ws.Range("J29").Value = 2000
$ws.Range("K29").Value = 2550
$ws.Range("L29").Value = 6000
$ws.Range("M29").Value = -2269
$ws.Range("N29").Value = -6562

# Row 38
$ws.Range("H38").Value = 42.833332
$ws.Range("I38").Value = 42.833332
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 128.499996
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 243.500004
$ws.Range("N38").ClearContents()

# Row 58
$ws.Range("H58").Value = 860.0909
$ws.Range("I58").Value = 275.625
$ws.Range("J58").Value = 1194.0714
$ws.Range("K58").Value = 826.875
$ws.Range("L58").Value = 3582.2142
$ws.Range("M58").Value = -676.875
$ws.Range("N58").Value = -3882.2142

# Row 112
$ws.Range("H112").Value = 6623.1816
$ws.Range("I112").Value = 633.3333
$ws.Range("J112").Value = 7568.9473
$ws.Range("K112").Value = 1899.9999
$ws.Range("L112").Value = 22706.8419
$ws.Range("M112").Value = -791.9999
$ws.Range("N112").Value = -24922.8419

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

# Row 64
$ws.Range("H64").Value = 2783
$ws.Range("I64").Value = 887.3333
$ws.Range("J64").Value = 3300
$ws.Range("K64").Value = 2661.9999
$ws.Range("L64").Value = 9900
$ws.Range("M64").Value = -2391.9999
$ws.Range("N64").Value = -10440

# Row 67
$ws.Range("H67").Value = 2783
$ws.Range("I67").Value = 887.3333
$ws.Range("J67").Value = 3300
$ws.Range("K67").Value = 2661.9999
$ws.Range("L67").Value = 9900
$ws.Range("M67").Value = -1725.9999
$ws.Range("N67").Value = -11772

# Row 140
$ws.Range("H140").Value = 1697.4445
$ws.Range("I140").Value = 1369.25
$ws.Range("J140").Value = 1960
$ws.Range("K140").Value = 4107.75
$ws.Range("L140").Value = 5880
$ws.Range("M140").Value = 1072.25
$ws.Range("N140").Value = -16240

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

# Row 43
$ws.Range("H43").Value = 5137.1113
$ws.Range("I43").Value = 808.5
$ws.Range("J43").Value = 8600
$ws.Range("K43").Value = 808.5
$ws.Range("L43").Value = 8600
$ws.Range("M43").Value = -657.5
$ws.Range("N43").Value = -8902

# Row 46
$ws.Range("H46").Value = 4000
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 4000
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 4000
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -4312

# Row 57
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

# Row 80
$ws.Range("H80").Value = 2679.2856
$ws.Range("I80").Value = 2451
$ws.Range("J80").Value = 3250
$ws.Range("K80").Value = 2451
$ws.Range("L80").Value = 3250
$ws.Range("M80").Value = -1453
$ws.Range("N80").Value = -5246

# Row 83
$ws.Range("H83").Value = 2679.2856
$ws.Range("I83").Value = 2451
$ws.Range("J83").Value = 3250
$ws.Range("K83").Value = 12255
$ws.Range("L83").Value = 16250
$ws.Range("M83").Value = -7263
$ws.Range("N83").Value = -26234

# Row 126
$ws.Range("H126").Value = 2165.5293
$ws.Range("I126").Value = 1400
$ws.Range("J126").Value = 2267.6
$ws.Range("K126").Value = 4200
$ws.Range("L126").Value = 6802.799999999999
$ws.Range("M126").Value = -1730
$ws.Range("N126").Value = -11742.8

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

# Row 7
$ws.Range("H7").Value = 1204599.5
$ws.Range("I7").Value = 2125.2
$ws.Range("J7").Value = 3208723.5
$ws.Range("K7").Value = 2125.2
$ws.Range("L7").Value = 3208723.5
$ws.Range("M7").Value = -2013.2
$ws.Range("N7").Value = -3208947.5

# Row 22
$ws.Range("H22").Value = 2000.2
$ws.Range("I22").Value = 2127.2727
$ws.Range("J22").Value = 1844.8889
$ws.Range("K22").Value = 2127.2727
$ws.Range("L22").Value = 1844.8889
$ws.Range("M22").Value = -1832.2727
$ws.Range("N22").Value = -2434.8889

# Row 27
$ws.Range("H27").Value = 2000.2
$ws.Range("I27").Value = 2127.2727
$ws.Range("J27").Value = 1844.8889
$ws.Range("K27").Value = 2127.2727
$ws.Range("L27").Value = 1844.8889
$ws.Range("M27").Value = -2020.2727
$ws.Range("N27").Value = -2058.8889

# Row 45
$ws.Range("H45").Value = 20000
$ws.Range("I45").Value = 10000
$ws.Range("K45").Value = 10000
$ws.Range("M45").Value = -9593

# Row 46
$ws.Range("H46").Value = 5609.1816
$ws.Range("I46").Value = 1131.6842
$ws.Range("J46").Value = 33966.668
$ws.Range("K46").Value = 1131.6842
$ws.Range("L46").Value = 33966.668
$ws.Range("M46").Value = -943.6841999999999
$ws.Range("N46").Value = -34342.668

# Row 48
$ws.Range("H48").Value = 30000
$ws.Range("J48").Value = 30000
$ws.Range("L48").Value = 30000
$ws.Range("N48").Value = -31322

# Row 55
$ws.Range("H55").Value = 4983
$ws.Range("I55").Value = 1057.0588
$ws.Range("J55").Value = 14517.429
$ws.Range("K55").Value = 1057.0588
$ws.Range("L55").Value = 14517.429
$ws.Range("M55").Value = -884.0588
$ws.Range("N55").Value = -14863.429

# Row 122
$ws.Range("H122").Value = 73532480
$ws.Range("I122").Value = 333335200
$ws.Range("J122").Value = 17860464
$ws.Range("K122").Value = 1000005600
$ws.Range("L122").Value = 53581392
$ws.Range("M122").Value = -1000003150
$ws.Range("N122").Value = -53586292

# Row 126
$ws.Range("H126").Value = 1204599.5
$ws.Range("I126").Value = 2125.2
$ws.Range("J126").Value = 3208723.5
$ws.Range("K126").Value = 6375.599999999999
$ws.Range("L126").Value = 9626170.5
$ws.Range("M126").Value = -3905.599999999999
$ws.Range("N126").Value = -9631110.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

# Row 81
$ws.Range("H81").Value = 20836260
$ws.Range("I81").Value = 1465.1666
$ws.Range("J81").Value = 41671056
$ws.Range("K81").Value = 2930.3332
$ws.Range("L81").Value = 83342112
$ws.Range("M81").Value = -1869.3332
$ws.Range("N81").Value = -83344234

# Row 84
$ws.Range("H84").Value = 20836260
$ws.Range("I84").Value = 1465.1666
$ws.Range("J84").Value = 41671056
$ws.Range("K84").Value = 14651.666
$ws.Range("L84").Value = 416710560
$ws.Range("M84").Value = -9347.666000000001
$ws.Range("N84").Value = -416721168

# Row 107
$ws.Range("H107").Value = 125000504
$ws.Range("I107").Value = 1000
$ws.Range("J107").Value = 166667000
$ws.Range("K107").Value = 1000
$ws.Range("L107").Value = 500001000
$ws.Range("M107").Value = -1080
$ws.Range("N107").Value = -500004840
